$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look numeric but must remain text (matches source inlineStr formatting)
$textCells = @('D4', 'D5', 'D6', 'D8', 'D10', 'D11', 'D12', 'D13', 'D14', 'D16', 'D19', 'D20', 'D21', 'D22', 'D23', 'D24', 'D25', 'D27', 'D28', 'D29', 'D30', 'D31', 'D33', 'D34', 'D35', 'D36', 'D37', 'D38', 'D39', 'D40', 'D43', 'D44', 'D45', 'D46', 'D48', 'D49', 'D50')
foreach ($c in $textCells) {
    $ws.Range($c).NumberFormat = "@"
}

# Apply updated values
$ws.Range('D2').Value = '67.887.12'
$ws.Range('E2').Value = '  -3.26%  '
$ws.Range('D3').Value = '3.332.07'
$ws.Range('E3').Value = '  -4.96%  '
$ws.Range('D4').Value = '0.996'
$ws.Range('E4').Value = '  -0.46%  '
$ws.Range('D5').Value = '603.62'
$ws.Range('E5').Value = '  +0.05%  '
$ws.Range('D6').Value = '163.42'
$ws.Range('E6').Value = '  -6.41%  '
$ws.Range('E7').Value = '  +0.05%  '
$ws.Range('D8').Value = '0.578'
$ws.Range('E8').Value = '  -5.28%  '
$ws.Range('D9').Value = '3.320.82'
$ws.Range('E9').Value = '  -5.13%  '
$ws.Range('D10').Value = '0.185'
$ws.Range('E10').Value = '  -4.22%  '
$ws.Range('D11').Value = '6.64'
$ws.Range('E11').Value = '  -8.08%  '
$ws.Range('D12').Value = '0.532'
$ws.Range('E12').Value = '  -8.65%  '
$ws.Range('D13').Value = '41.73'
$ws.Range('E13').Value = '  -9.47%  '
$ws.Range('D14').Value = '0.0000256'
$ws.Range('E14').Value = '  -6.75%  '
$ws.Range('D15').Value = '3.869.42'
$ws.Range('E15').Value = '  -5.09%  '
$ws.Range('B16').Value = 'Polkadot'
$ws.Range('C16').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D16').Value = '7.76'
$ws.Range('E16').Value = '  -6.23%  '
$ws.Range('D17').Value = '67.792.30'
$ws.Range('E17').Value = '  -3.60%  '
$ws.Range('B18').Value = 'WrappedEther'
$ws.Range('C18').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D18').Value = '3.325.64'
$ws.Range('E18').Value = '  -5.43%  '
$ws.Range('B19').Value = 'BitcoinCash'
$ws.Range('C19').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D19').Value = '558.40'
$ws.Range('E19').Value = '  -8.45%  '
$ws.Range('B20').Value = 'TRON'
$ws.Range('C20').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D20').Value = '0.118'
$ws.Range('E20').Value = '  -1.31%  '
$ws.Range('D21').Value = '16.24'
$ws.Range('E21').Value = '  -6.03%  '
$ws.Range('D22').Value = '0.804'
$ws.Range('E22').Value = '  -7.80%  '
$ws.Range('D23').Value = '8.55'
$ws.Range('E23').Value = '  -5.33%  '
$ws.Range('D24').Value = '89.88'
$ws.Range('E24').Value = '  -7.76%  '
$ws.Range('D25').Value = '14.25'
$ws.Range('E25').Value = '  -8.71%  '
$ws.Range('E26').Value = '  -6.62%  '
$ws.Range('D27').Value = '0.999'
$ws.Range('E27').Value = '  -0.05%  '
$ws.Range('D28').Value = '31.25'
$ws.Range('E28').Value = '  -8.04%  '
$ws.Range('D29').Value = '2.28'
$ws.Range('E29').Value = '  -10.48%  '
$ws.Range('D30').Value = '8.14'
$ws.Range('E30').Value = '  -8.95%  '
$ws.Range('D31').Value = '7.35'
$ws.Range('E31').Value = '  -9.18%  '
$ws.Range('E32').Value = '  -7.68%  '
$ws.Range('D33').Value = '2.65'
$ws.Range('E33').Value = '  -10.99%  '
$ws.Range('D34').Value = '575.38'
$ws.Range('E34').Value = '  -9.29%  '
$ws.Range('D35').Value = '6.24'
$ws.Range('E35').Value = '  -9.18%  '
$ws.Range('D36').Value = '1.00'
$ws.Range('E36').Value = '  +0.07%  '
$ws.Range('D37').Value = '54.97'
$ws.Range('E37').Value = '  -2.93%  '
$ws.Range('D38').Value = '9.95'
$ws.Range('E38').Value = '  -6.89%  '
$ws.Range('D39').Value = '0.0909'
$ws.Range('E39').Value = '  -8.51%  '
$ws.Range('D40').Value = '0.0448'
$ws.Range('E40').Value = '  -4.91%  '
$ws.Range('E41').Value = '  -3.72%  '
$ws.Range('D42').Value = '3.086.81'
$ws.Range('E42').Value = '  -7.91%  '
$ws.Range('D43').Value = '2.90'
$ws.Range('E43').Value = '  -18.80%  '
$ws.Range('D44').Value = '2.65'
$ws.Range('E44').Value = '  -8.46%  '
$ws.Range('D45').Value = '29.15'
$ws.Range('E45').Value = '  -9.13%  '
$ws.Range('D46').Value = '0.278'
$ws.Range('E46').Value = '  -10.10%  '
$ws.Range('D47').Value = '0.0₃0623'
$ws.Range('E47').Value = '  -15.47%  '
$ws.Range('D48').Value = '2.25'
$ws.Range('E48').Value = '  -11.50%  '
$ws.Range('D49').Value = '0.121'
$ws.Range('E49').Value = '  -6.61%  '
$ws.Range('D50').Value = '129.35'
$ws.Range('E50').Value = '  -3.38%  '
$ws.Range('E51').Value = '  -0.02%  '
